# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-19 14:14:52
#
# The Biochemistry Lab/CBL session-1 row (row 7) got its attendance recorded,
# which ripples into the "Class Statistics" summary block (L6:L10) and the
# "Group Statistics" table (row 15). A couple of "Recorded By" cells also got
# their email lists re-ordered (reporting refresh), and G15 / G2 / G3 / G4 do
# the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recorded-By email lists: reordered (same people, refreshed sort) -------
$ws.Range("G2").Value = "System, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 7 (BIOCHEMISTRY LAB/CBL, session 1) just got recorded --------------
# Pick up the "Recorded" look (green fill) by copying the format from a row
# that already carries it (row 2), then overwrite the now-populated cells.
$ws.Range("A2:I2").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)
$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("H7").Value = "19/251"
$ws.Range("I7").Value = "Recorded"

# --- Class Statistics block (K6:L10) ----------------------------------------
# L9/L10 hold their percentages as literal text (not real % numbers), so force
# a text number-format before assigning or Excel will coerce "34.5%" into the
# number 0.345.
$ws.Range("L6").Value = 10        # Recorded Sessions:  9 -> 10
$ws.Range("L7").Value = 1         # Missing Sessions:   2 -> 1
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "34.5%"   # Coverage %:         31.0% -> 34.5%
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "23.7%"  # Average Attendance %: 25.5% -> 23.7%

# --- Group Statistics table (row 15) mirrors the same recalculation --------
$ws.Range("O15").Value = 10
$ws.Range("P15").Value = 1
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "34.5%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "23.7%"
